# Calendario (Interfaz Reportes 100%)
# Row 13 = "ReportesController" -> progress for weeks 4 & 5 (E13, F13) reaches 100%.
# Also, the per-column "traffic-light" (cellIs) conditional formatting that used
# to cover B2:E14 and, separately, F2:F14 is unified into a single B2:F14 rule
# set (matching the diff's merged <conditionalFormatting sqref="B2:F14"> block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Data edit: ReportesController reaches 100% for the last two weeks ---
$ws.Range("E13").Value = 100
$ws.Range("F13").Value = 100

# --- 2. Conditional formatting: merge B2:E14 + F2:F14 "cellIs" rules into B2:F14 ---

# 2a. Remove the old cellIs ("Between") rules that applied only to B2:E14,
#     leaving the colorScale rules (and the F2:F14 cellIs rules) untouched.
$rngBE = $ws.Range("B2:E14")
$more = $true
while ($more) {
    $more = $false
    $count = $rngBE.FormatConditions.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc = $rngBE.FormatConditions.Item($i)
        if ($fc.Type -eq 1 -and $fc.AppliesTo.Address() -eq "`$B`$2:`$E`$14") {
            $fc.Delete()
            $more = $true
            break
        }
    }
}

# 2b. Extend the existing F2:F14 cellIs rules to cover B2:F14 instead, so the
#     same three rules (same colors/dxf, same priority order) now apply to the
#     whole B:F block, exactly like the merged rule set in the target file.
$rngF = $ws.Range("F2:F14")
$newRange = $ws.Range("B2:F14")
for ($i = 1; $i -le $rngF.FormatConditions.Count; $i++) {
    $fc = $rngF.FormatConditions.Item($i)
    if ($fc.Type -eq 1) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- 3. Leave the selection on F12, matching the saved cursor position ---
$ws.Range("F12").Select()
